$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A3").Value = "Test Case ID"
Write-Host "ok"
